# Business Category List.xlsx edit
# Commit: "add real estate, agency, job, consulting, startup, ecommerce, startup"
#
# The tracker rows for Landing Page, Real Estate, Agency, Job, Consulting,
# Startup and Ecommerce (rows 39-45) are marked as completed, assigned to
# "Fahri" (column D) with status "DONE" (column E) - matching the pattern
# already used by every other finished row (e.g. row 46 "Industrial").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 39-45 (categories: Landing Page, Real Estate, Agency, Job,
# Consulting, Startup, Ecommerce) as done.
$ws.Range("D39:D45").Value = "Fahri"
$ws.Range("E39:E45").Value = "DONE"

# Reflect the author's scroll/zoom/selection state at save time.
$ws.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 105
$ws.Range("E38:E45").Select()
